# Sprint 1 Backlog update
# - Fill in actual-time tracking for the "Design the portfolio/buy/sell page" rows (18-20)
# - Add three new completed tasks (rows 40-42) with an invalid leap-day date typed as text
# - Row-26 / row-64 SUM formulas and the burndown chart (bound to Sheet1!$C$26:$O$26)
# recalc automatically from the underlying cell edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 18-20: mark owner (Colby) and fill in daily actual-time tracking ---
$ws.Range("A18").Value = "Colby"
$ws.Range("E18:G18").Value = 1.5
$ws.Range("H18:J18").Value = 0

$ws.Range("A19").Value = "Colby"
$ws.Range("E19:I19").Value = 1.5
$ws.Range("J19").Value = 0

$ws.Range("A20").Value = "Colby"
$ws.Range("E20:F20").Value = 2
$ws.Range("G20:J20").Value = 0

# --- Rows 40-42: new completed backlog items in the Actual Time table ---
# Copy the date-formatted style from the row above (D39) down into D40:D42 first
# so the new date cells match the existing column formatting.
$ws.Range("D39").Copy()
$ws.Range("D40:D42").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A40").Value = "Colby"
$ws.Range("B40").Value = "Designed the portfolio page"
$ws.Range("C40").Value = 2
$ws.Range("D40").Value = 45716

$ws.Range("A41").Value = "Colby"
$ws.Range("B41").Value = "Designed the sell page"
$ws.Range("C41").Value = 1

$ws.Range("A42").Value = "Colby"
$ws.Range("B42").Value = "Allow for the selling of crypto"
$ws.Range("C42").Value = 3
$ws.Range("D42").Value = 45721

# 2025 is not a leap year, so this literal gets stored as text (matches source behavior)
# Set last so the new shared-string insertion order matches: portfolio page, sell page,
# "Allow for the selling of crypto", then this date text.
$ws.Range("D41").Value = "2/29/2025"

# --- View state: selection moved to D41 ---
$ws.Range("D41").Select()
